$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$arr = New-Object 'object[,]' 16,20
$arr[0,0] = "FAPs"
$arr[0,1] = "Hgf"
$arr[0,2] = "Cd44"
$arr[0,3] = "ECs"
$arr[0,4] = 2
$arr[0,5] = 0.6666666666666666
$arr[0,6] = 0.1088396666666667
$arr[0,7] = 0.326519
$arr[0,8] = 0.002750770615347974
$arr[0,9] = 0.002750770615347974
$arr[0,10] = 3
$arr[0,11] = 1
$arr[0,12] = 31.82741333333333
$arr[0,13] = 95.48223999999999
$arr[0,14] = 0.114390792932228
$arr[0,15] = 0.114390792932228
$arr[0,16] = 3.464085058062222
$arr[0,17] = 31.17676552256
$arr[0,18] = 0.0003146628318643275
$arr[0,19] = 0.0003146628318643276
$arr[1,0] = "FAPs"
$arr[1,1] = "Hgf"
$arr[1,2] = "Cd44"
$arr[1,3] = "FAPs"
$arr[1,4] = 2
$arr[1,5] = 0.6666666666666666
$arr[1,6] = 0.1088396666666667
$arr[1,7] = 0.326519
$arr[1,8] = 0.002750770615347974
$arr[1,9] = 0.002750770615347974
$arr[1,10] = 3
$arr[1,11] = 1
$arr[1,12] = 85.46317833333335
$arr[1,13] = 256.389535
$arr[1,14] = 0.307162904935779
$arr[1,15] = 0.307162904935779
$arr[1,16] = 9.30178384207389
$arr[1,17] = 83.71605457866501
$arr[1,18] = 0.000844934693022264
$arr[1,19] = 0.0008449346930222641
$arr[2,0] = "FAPs"
$arr[2,1] = "Hgf"
$arr[2,2] = "Cd44"
$arr[2,3] = "M2"
$arr[2,4] = 2
$arr[2,5] = 0.6666666666666666
$arr[2,6] = 0.1088396666666667
$arr[2,7] = 0.326519
$arr[2,8] = 0.002750770615347974
$arr[2,9] = 0.002750770615347974
$arr[2,10] = 3
$arr[2,11] = 1
$arr[2,12] = 122.2478306666667
$arr[2,13] = 366.743492
$arr[2,14] = 0.4393704929064738
$arr[2,15] = 0.4393704929064738
$arr[2,16] = 13.30541314048311
$arr[2,17] = 119.748718264348
$arr[2,18] = 0.001208607441138083
$arr[2,19] = 0.001208607441138084
$arr[3,0] = "FAPs"
$arr[3,1] = "Hgf"
$arr[3,2] = "Cd44"
$arr[3,3] = "sCs"
$arr[3,4] = 2
$arr[3,5] = 0.6666666666666666
$arr[3,6] = 0.1088396666666667
$arr[3,7] = 0.326519
$arr[3,8] = 0.002750770615347974
$arr[3,9] = 0.002750770615347974
$arr[3,10] = 3
$arr[3,11] = 1
$arr[3,12] = 38.69562533333333
$arr[3,13] = 116.086876
$arr[3,14] = 0.1390758092255191
$arr[3,15] = 0.1390758092255191
$arr[3,16] = 4.211618962738222
$arr[3,17] = 37.904570664644
$arr[3,18] = 0.0003825656493232986
$arr[3,19] = 0.0003825656493232987
$arr[4,0] = "M2"
$arr[4,1] = "Hgf"
$arr[4,2] = "Cd44"
$arr[4,3] = "ECs"
$arr[4,4] = 3
$arr[4,5] = 1
$arr[4,6] = 15.07419333333333
$arr[4,7] = 45.22258
$arr[4,8] = 0.3809791902285103
$arr[4,9] = 0.3809791902285103
$arr[4,10] = 3
$arr[4,11] = 1
$arr[4,12] = 31.82741333333333
$arr[4,13] = 95.48223999999999
$arr[4,14] = 0.114390792932228
$arr[4,15] = 0.114390792932228
$arr[4,16] = 479.7725818865777
$arr[4,17] = 4317.9532369792
$arr[4,18] = 0.04358051166091743
$arr[4,19] = 0.04358051166091744
$arr[5,0] = "M2"
$arr[5,1] = "Hgf"
$arr[5,2] = "Cd44"
$arr[5,3] = "FAPs"
$arr[5,4] = 3
$arr[5,5] = 1
$arr[5,6] = 15.07419333333333
$arr[5,7] = 45.22258
$arr[5,8] = 0.3809791902285103
$arr[5,9] = 0.3809791902285103
$arr[5,10] = 3
$arr[5,11] = 1
$arr[5,12] = 85.46317833333335
$arr[5,13] = 256.389535
$arr[5,14] = 0.307162904935779
$arr[5,15] = 0.307162904935779
$arr[5,16] = 1288.288473077811
$arr[5,17] = 11594.5962577003
$arr[5,18] = 0.11702267479067
$arr[5,19] = 0.11702267479067
$arr[6,0] = "M2"
$arr[6,1] = "Hgf"
$arr[6,2] = "Cd44"
$arr[6,3] = "M2"
$arr[6,4] = 3
$arr[6,5] = 1
$arr[6,6] = 15.07419333333333
$arr[6,7] = 45.22258
$arr[6,8] = 0.3809791902285103
$arr[6,9] = 0.3809791902285103
$arr[6,10] = 3
$arr[6,11] = 1
$arr[6,12] = 122.2478306666667
$arr[6,13] = 366.743492
$arr[6,14] = 0.4393704929064738
$arr[6,15] = 0.4393704929064738
$arr[6,16] = 1842.787434049929
$arr[6,17] = 16585.08690644936
$arr[6,18] = 0.1673910145978098
$arr[6,19] = 0.1673910145978099
$arr[7,0] = "M2"
$arr[7,1] = "Hgf"
$arr[7,2] = "Cd44"
$arr[7,3] = "sCs"
$arr[7,4] = 3
$arr[7,5] = 1
$arr[7,6] = 15.07419333333333
$arr[7,7] = 45.22258
$arr[7,8] = 0.3809791902285103
$arr[7,9] = 0.3809791902285103
$arr[7,10] = 3
$arr[7,11] = 1
$arr[7,12] = 38.69562533333333
$arr[7,13] = 116.086876
$arr[7,14] = 0.1390758092255191
$arr[7,15] = 0.1390758092255191
$arr[7,16] = 583.3053374288978
$arr[7,17] = 5249.74803686008
$arr[7,18] = 0.05298498917911307
$arr[7,19] = 0.05298498917911307
$arr[8,0] = "sCs"
$arr[8,1] = "Hgf"
$arr[8,2] = "Cd44"
$arr[8,3] = "ECs"
$arr[8,4] = 3
$arr[8,5] = 1
$arr[8,6] = 23.69325166666667
$arr[8,7] = 71.079755
$arr[8,8] = 0.5988138558556568
$arr[8,9] = 0.5988138558556569
$arr[8,10] = 3
$arr[8,11] = 1
$arr[8,12] = 31.82741333333333
$arr[8,13] = 95.48223999999999
$arr[8,14] = 0.114390792932228
$arr[8,15] = 0.114390792932228
$arr[8,16] = 754.0949140056888
$arr[8,17] = 6786.854226051199
$arr[8,18] = 0.06849879179013348
$arr[8,19] = 0.0684987917901335
$arr[9,0] = "sCs"
$arr[9,1] = "Hgf"
$arr[9,2] = "Cd44"
$arr[9,3] = "FAPs"
$arr[9,4] = 3
$arr[9,5] = 1
$arr[9,6] = 23.69325166666667
$arr[9,7] = 71.079755
$arr[9,8] = 0.5988138558556568
$arr[9,9] = 0.5988138558556569
$arr[9,10] = 3
$arr[9,11] = 1
$arr[9,12] = 85.46317833333335
$arr[9,13] = 256.389535
$arr[9,14] = 0.307162904935779
$arr[9,15] = 0.307162904935779
$arr[9,16] = 2024.900592484881
$arr[9,17] = 18224.10533236393
$arr[9,18] = 0.1839334034804184
$arr[9,19] = 0.1839334034804184
$arr[10,0] = "sCs"
$arr[10,1] = "Hgf"
$arr[10,2] = "Cd44"
$arr[10,3] = "M2"
$arr[10,4] = 3
$arr[10,5] = 1
$arr[10,6] = 23.69325166666667
$arr[10,7] = 71.079755
$arr[10,8] = 0.5988138558556568
$arr[10,9] = 0.5988138558556569
$arr[10,10] = 3
$arr[10,11] = 1
$arr[10,12] = 122.2478306666667
$arr[10,13] = 366.743492
$arr[10,14] = 0.4393704929064738
$arr[10,15] = 0.4393704929064738
$arr[10,16] = 2896.448617689385
$arr[10,17] = 26068.03755920446
$arr[10,18] = 0.2631011390065261
$arr[10,19] = 0.2631011390065262
$arr[11,0] = "sCs"
$arr[11,1] = "Hgf"
$arr[11,2] = "Cd44"
$arr[11,3] = "sCs"
$arr[11,4] = 3
$arr[11,5] = 1
$arr[11,6] = 23.69325166666667
$arr[11,7] = 71.079755
$arr[11,8] = 0.5988138558556568
$arr[11,9] = 0.5988138558556569
$arr[11,10] = 3
$arr[11,11] = 1
$arr[11,12] = 38.69562533333333
$arr[11,13] = 116.086876
$arr[11,14] = 0.1390758092255191
$arr[11,15] = 0.1390758092255191
$arr[11,16] = 916.825189421709
$arr[11,17] = 8251.42670479538
$arr[11,18] = 0.08328052157857883
$arr[11,19] = 0.08328052157857885
$arr[12,0] = "ECs"
$arr[12,1] = "Hgf"
$arr[12,2] = "Cd44"
$arr[12,3] = "ECs"
$arr[12,4] = 3
$arr[12,5] = 1
$arr[12,6] = 0.6906883333333335
$arr[12,7] = 2.072065
$arr[12,8] = 0.01745618330048481
$arr[12,9] = 0.01745618330048481
$arr[12,10] = 3
$arr[12,11] = 1
$arr[12,12] = 31.82741333333333
$arr[12,13] = 95.48223999999999
$arr[12,14] = 0.114390792932228
$arr[12,15] = 0.114390792932228
$arr[12,16] = 21.98282306951111
$arr[12,17] = 197.8454076256
$arr[12,18] = 0.001996826649312774
$arr[12,19] = 0.001996826649312775
$arr[13,0] = "ECs"
$arr[13,1] = "Hgf"
$arr[13,2] = "Cd44"
$arr[13,3] = "FAPs"
$arr[13,4] = 3
$arr[13,5] = 1
$arr[13,6] = 0.6906883333333335
$arr[13,7] = 2.072065
$arr[13,8] = 0.01745618330048481
$arr[13,9] = 0.01745618330048481
$arr[13,10] = 3
$arr[13,11] = 1
$arr[13,12] = 85.46317833333335
$arr[13,13] = 256.389535
$arr[13,14] = 0.307162904935779
$arr[13,15] = 0.307162904935779
$arr[13,16] = 59.02842020441946
$arr[13,17] = 531.2557818397751
$arr[13,18] = 0.00536189197166835
$arr[13,19] = 0.00536189197166835
$arr[14,0] = "ECs"
$arr[14,1] = "Hgf"
$arr[14,2] = "Cd44"
$arr[14,3] = "M2"
$arr[14,4] = 3
$arr[14,5] = 1
$arr[14,6] = 0.6906883333333335
$arr[14,7] = 2.072065
$arr[14,8] = 0.01745618330048481
$arr[14,9] = 0.01745618330048481
$arr[14,10] = 3
$arr[14,11] = 1
$arr[14,12] = 122.2478306666667
$arr[14,13] = 366.743492
$arr[14,14] = 0.4393704929064738
$arr[14,15] = 0.4393704929064738
$arr[14,16] = 84.43515041677557
$arr[14,17] = 759.9163537509801
$arr[14,18] = 0.007669731860999768
$arr[14,19] = 0.007669731860999769
$arr[15,0] = "ECs"
$arr[15,1] = "Hgf"
$arr[15,2] = "Cd44"
$arr[15,3] = "sCs"
$arr[15,4] = 3
$arr[15,5] = 1
$arr[15,6] = 0.6906883333333335
$arr[15,7] = 2.072065
$arr[15,8] = 0.01745618330048481
$arr[15,9] = 0.01745618330048481
$arr[15,10] = 3
$arr[15,11] = 1
$arr[15,12] = 38.69562533333333
$arr[15,13] = 116.086876
$arr[15,14] = 0.1390758092255191
$arr[15,15] = 0.1390758092255191
$arr[15,16] = 26.72661696877111
$arr[15,17] = 240.53955271894
$arr[15,18] = 0.002427732818503918
$arr[15,19] = 0.002427732818503918

$ws.Range("A2:T17").Value = $arr
